$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the base reference values used throughout the sheet's formulas
$ws.Range("F1").Value = 53
$ws.Range("G1").Value = 1.9

# Force a full recalculation so all dependent formulas (columns E-J, row 24,
# and the chart's cached series) reflect the new F1/G1 values
$excel.CalculateFullRebuild()

# Update the view state: scroll so column E is the left-most visible column,
# and move/extend the active selection
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("H3:J21").Select()
$excel.ActiveCell = $ws.Range("J21")
